$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("E3").Value = "-"
$ws.Range("F3").Value = "-"

# Row 4
$ws.Range("C4").Value = "Mecanica material"
$ws.Range("E4").Value = "CAD"
$ws.Range("F4").Value = "MTRM"

# Row 6
$ws.Range("C6").Value = "Mecanica material"
$ws.Range("F6").Value = "MTRM"
